$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 701.09
$ws.Range("C3").Value = 724.25
$ws.Range("C4").Value = 700.77
$ws.Range("C5").Value = 719.63
$ws.Range("C6").Value = 719.63
